$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.111543774604797
$ws.Range("B1").Value = 2.506280183792114
$ws.Range("C1").Value = 6.215395927429199
$ws.Range("D1").Value = 2.181187152862549
$ws.Range("E1").Value = 1.256453990936279
